# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2404"
#   "<name>_new" -> "<name>_FV2410"
# Then turn the data range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21

# 1. Rename the header cells in row 1.
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2404'
        $newVal = $newVal -replace '_new$', '_FV2410'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# 2. Convert the header/data range into a native Excel Table ("Table1").
$tableRange = $ws.Range("A1:U60")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3. Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Activate()
$ws.Rows(2).Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
